$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.943.26"
$ws.Range("D2").ClearFormats()
$ws.Range("D3").Value = "'3.209.86"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'  +0.94%  "
$ws.Range("E3").ClearFormats()
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "'  -0.01%  "
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'602.97"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'  +4.70%  "
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = "'151.73"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'  +0.34%  "
$ws.Range("E6").ClearFormats()
$ws.Range("D8").Value = "'3.207.51"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'  +0.97%  "
$ws.Range("E8").ClearFormats()
$ws.Range("D9").Value = "'0.537"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'  +1.65%  "
$ws.Range("E9").ClearFormats()
$ws.Range("E10").Value = "'  -1.62%  "
$ws.Range("E10").ClearFormats()
$ws.Range("E11").Value = "'  -1.20%  "
$ws.Range("E11").ClearFormats()
$ws.Range("D12").Value = "'0.511"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'  +0.51%  "
$ws.Range("E12").ClearFormats()
$ws.Range("D13").Value = "'0.0000271"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'  -2.02%  "
$ws.Range("E13").ClearFormats()
$ws.Range("D14").Value = "'38.62"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'  +1.14%  "
$ws.Range("E14").ClearFormats()
$ws.Range("D15").Value = "'3.738.69"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'  +1.07%  "
$ws.Range("E15").ClearFormats()
$ws.Range("D16").Value = "'66.026.38"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'  +1.21%  "
$ws.Range("E16").ClearFormats()
$ws.Range("D17").Value = "'7.40"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'  +2.62%  "
$ws.Range("E17").ClearFormats()
$ws.Range("D18").Value = "'3.210.82"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'  +0.75%  "
$ws.Range("E18").ClearFormats()
$ws.Range("E19").Value = "'  +0.42%  "
$ws.Range("E19").ClearFormats()
$ws.Range("D20").Value = "'512.34"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "'  -0.44%  "
$ws.Range("E20").ClearFormats()
$ws.Range("D21").Value = "'15.77"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "'  +5.51%  "
$ws.Range("E21").ClearFormats()
$ws.Range("D22").Value = "'0.737"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'  +0.25%  "
$ws.Range("E22").ClearFormats()
$ws.Range("D23").Value = "'15.20"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'  -1.98%  "
$ws.Range("E23").ClearFormats()
$ws.Range("E24").Value = "'  +1.50%  "
$ws.Range("E24").ClearFormats()
$ws.Range("D25").Value = "'85.39"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "'  +0.26%  "
$ws.Range("E25").ClearFormats()
$ws.Range("D27").Value = "'9.28"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "'  +2.29%  "
$ws.Range("E27").ClearFormats()
$ws.Range("D28").Value = "'3.02"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "'  +3.58%  "
$ws.Range("E28").ClearFormats()
$ws.Range("D29").Value = "'2.24"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "'  +2.08%  "
$ws.Range("E29").ClearFormats()
$ws.Range("D30").Value = "'2.88"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "'  +3.81%  "
$ws.Range("E30").ClearFormats()
$ws.Range("D31").Value = "'6.81"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "'  +7.50%  "
$ws.Range("E31").ClearFormats()
$ws.Range("D32").Value = "'28.17"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "'  +0.04%  "
$ws.Range("E32").ClearFormats()
$ws.Range("E33").Value = "'  +0.17%  "
$ws.Range("E33").ClearFormats()
$ws.Range("E34").Value = "'  +0.10%  "
$ws.Range("E34").ClearFormats()
$ws.Range("E35").Value = "'  -0.97%  "
$ws.Range("E35").ClearFormats()
$ws.Range("D36").Value = "'55.48"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "'  -0.60%  "
$ws.Range("E36").ClearFormats()
$ws.Range("D37").Value = "'0.0923"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "'  +3.06%  "
$ws.Range("E37").ClearFormats()
$ws.Range("D38").Value = "'488.05"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "'  +1.52%  "
$ws.Range("E38").ClearFormats()
$ws.Range("D39").Value = "'0.0424"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "'  +0.09%  "
$ws.Range("E39").ClearFormats()
$ws.Range("E40").Value = "'  -4.01%  "
$ws.Range("E40").ClearFormats()
$ws.Range("D41").Value = "'8.89"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'  +2.74%  "
$ws.Range("E41").ClearFormats()
$ws.Range("D42").Value = "'3.029.43"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'  -3.50%  "
$ws.Range("E42").ClearFormats()
$ws.Range("E43").Value = "'  -0.35%  "
$ws.Range("E43").ClearFormats()
$ws.Range("D44").Value = "'0.294"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'  +1.67%  "
$ws.Range("E44").ClearFormats()
$ws.Range("D45").Value = "'0.0₃0648"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'  +6.17%  "
$ws.Range("E45").ClearFormats()
$ws.Range("D46").Value = "'2.46"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'  -0.59%  "
$ws.Range("E46").ClearFormats()
$ws.Range("D47").Value = "'29.16"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "'  -0.97%  "
$ws.Range("E47").ClearFormats()
$ws.Range("E48").Value = "'  +0.04%  "
$ws.Range("E48").ClearFormats()
$ws.Range("D49").Value = "'0.116"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "'  +0.44%  "
$ws.Range("E49").ClearFormats()
$ws.Range("E50").Value = "'  +1.38%  "
$ws.Range("E50").ClearFormats()
$ws.Range("D51").Value = "'120.46"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "'  -1.39%  "
$ws.Range("E51").ClearFormats()
